$d = $word.ActiveDocument
$d.Content.Find.Execute("ตาราง …", $true, $false, $false, $false, $false, $true, 1, $false, "ตาราง ที่ 1", 2)
